$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '62.513.88'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -2.91%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.379.54'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -3.72%  '

$ws.Range("E4").Value = '  -0.09%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '572.53'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -3.30%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '125.42'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -6.81%  '

$ws.Range("E7").Value = '  +0.01%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.379.89'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -3.71%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.476'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.48%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.28'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -4.61%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.376'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -3.27%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.951.39'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.83%  '

$ws.Range("E14").Value = '  -0.97%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.371.17'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -3.94%  '

$ws.Range("E16").Value = '  -6.03%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '62.513.00'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.91%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '24.37'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -5.35%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '9.24'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -7.51%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.63'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.25%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.08'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -4.20%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '372.87'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -5.45%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.554'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -4.44%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.512.64'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -3.75%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.00'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.11%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '71.60'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -4.08%  '

$ws.Range("E27").Value = '  -10.22%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.998'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.09%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.92'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -6.68%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.11'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -7.12%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.78'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -5.80%  '

$ws.Range("E32").Value = '  -0.02%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.38'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -6.99%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.410.33'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -3.62%  '

$ws.Range("E35").Value = '  -6.40%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '22.66'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.98%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.24'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.46%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '166.66'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.16%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.64'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -4.58%  '

$ws.Range("E40").Value = '  -5.64%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0752'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -4.56%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.998'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.21%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.764'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -5.89%  '

$ws.Range("E44").Value = '  -1.29%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.23'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -4.98%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '22.64'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -9.56%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.54'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -7.23%  '

$ws.Range("E48").Value = '  -8.36%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '6.58'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.28%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.239.86'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -5.79%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.838'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -8.05%  '
